$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New donor rows to append (rows 12-14).
$newRows = @(
    @{
        A = "69234632a9d7b28957d0eb0c"
        B = "mno"
        C = "mno@gmail.com"
        D = "9481824919"
        E = "A+"
        F = "Mangalore"
        G = 0
        H = 0
        I = $false
        J = $false
        K = $true
        L = 45984.96320155093
    },
    @{
        A = "69234680a9d7b28957d0eb11"
        B = "Me"
        C = "23a43.bhavish@sjec.ac.in"
        D = "8904534919"
        E = "A+"
        F = "Mangalore"
        G = 0
        H = 0
        I = $false
        J = $false
        K = $true
        L = 45984.9641012037
    },
    @{
        A = "6923b7181c52f11af60577df"
        B = "mailtrap"
        C = "lifelink@system.com"
        D = "1234456789"
        E = "A+"
        F = "Mangalore"
        G = 0
        H = 0
        I = $false
        J = $false
        K = $true
        L = 45985.297717719906
    }
)

$row = 12
foreach ($rec in $newRows) {
    $ws.Cells.Item($row, 1).Value = $rec.A
    $ws.Cells.Item($row, 2).Value = $rec.B
    $ws.Cells.Item($row, 3).Value = $rec.C

    # Phone numbers must stay text (they look numeric). Force Text format
    # before entry so the stored value is a string, then drop the explicit
    # style again so the cell matches the other plain text cells (no `s`).
    $dcell = $ws.Cells.Item($row, 4)
    $dcell.NumberFormat = "@"
    $dcell.Value = $rec.D
    $dcell.Style = "Normal"

    $ws.Cells.Item($row, 5).Value = $rec.E
    $ws.Cells.Item($row, 6).Value = $rec.F
    $ws.Cells.Item($row, 7).Value = $rec.G
    $ws.Cells.Item($row, 8).Value = $rec.H
    $ws.Cells.Item($row, 9).Value = $rec.I
    $ws.Cells.Item($row, 10).Value = $rec.J
    $ws.Cells.Item($row, 11).Value = $rec.K
    $ws.Cells.Item($row, 12).Value = $rec.L
    $row++
}

# The final new row (14) inherits the date number format that used to live
# on L11 (copy the format across so the same style index is reused), then
# L11 reverts to the default/general style.
$ws.Cells.Item(11, 12).Copy()
$ws.Cells.Item(14, 12).PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item(14, 12).Value = 45985.297717719906

$ws.Cells.Item(11, 12).Style = "Normal"
